$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("統計")

$newRow = 6

$ws.Cells.Item($newRow, 1).Value = "2025-08-27T18:32:07.930318"
$ws.Cells.Item($newRow, 2).Value = 12
$ws.Cells.Item($newRow, 3).Value = "全案件リスト"
$ws.Cells.Item($newRow, 4).Value = 75
$ws.Cells.Item($newRow, 5).Value = 3
$ws.Cells.Item($newRow, 6).Value = 6
$ws.Cells.Item($newRow, 7).Value = 12
